# "Make the total right according to the text files"
#
# Sheet1!G2 held 42 but should be 42.4. Sheet1!F2 (=F1+G2) is a formula,
# so updating G2 alone causes Excel to recompute F2 (92 -> 92.4) on
# recalculation - matching the target diff for both cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("G2").Value = 42.4
